$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 "paragraph" and B2 "Hello World" stay the same.

# New cell B12 "Let's try this out"
$ws.Range("B12").Value = "Let's try this out"

# New cell B3 "hehe"
$ws.Range("B3").Value = "hehe"

# A8 "// Ignore" -> leading spaces preserved: "     // Ignore"
$ws.Range("A8").Value = "     // Ignore"

# A9 "p" -> leading spaces preserved: "     p"
$ws.Range("A9").Value = "     p"

# B9 stays 18 (unchanged)

# C9 "purple" moves to D9
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = "purple"

# B10 "Hello Again" stays the same.

# Update selection to match the recorded state (D9)
$ws.Range("D9").Select()
